# Workbook "Fruta, Vega Modelo de Temuco - Chirimoya"
# Commit: "Fruta / hortaliza, semanal" - weekly update that inserts two new
# daily price records at the top of the Chirimoya dataset (rows 104-105),
# pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 104-105; this shifts the former rows
# 104..205 down to 106..207 (matches dimension change T205 -> T207).
$ws.Rows("104:105").Insert()

# After the insert, row 106 holds what used to be row 104's data, and
# row 107 holds what used to be row 105's data. Duplicate those rows back
# up into the newly inserted 104/105 rows as a starting point, then apply
# the specific field updates for the new records.
$ws.Range("A106:T106").Copy()
$ws.Range("A104:T104").PasteSpecial()

$ws.Range("A107:T107").Copy()
$ws.Range("A105:T105").PasteSpecial()

# New record in row 104: same as the old row 104 but with the newer date.
$ws.Range("D104").Value = 45175

# New record in row 105: same as the old row 105 but with the newer date
# and an updated Volumen (M) value.
$ws.Range("D105").Value = 45175
$ws.Range("M105").Value = 90
